# Power Storage Investment Option
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# Add the new "Power_storage" row of data (row 27) to the Investment_Cost table
$ws.Range("A27").Value = "Power_storage"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1

# Reflect the new active cell selection on the sheet
$ws.Activate()
$ws.Range("B28").Select()
